$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 05:17"

# --- Pakistan (row 26) ---
$ws.Range("B26").Value = 325480
$ws.Range("C26").Value = 736
$ws.Range("D26").Value = 309136
$ws.Range("E26").Value = 9642
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = 6702

# --- Row 29/30: Belgica & Paises Bajos swap places, with refreshed figures ---
# Row 29 now shows Belgica's (updated) data …
$ws.Range("A29").Value = "Belgica"
$ws.Range("B29").Value = 253386
$ws.Range("C29").Value = 13227
$ws.Range("D29").Value = 21717
$ws.Range("E29").Value = 221130
$ws.Range("G29").Value = 50
$ws.Range("H29").Value = 10539

# … and row 30 now shows Paises Bajos (with the figures Belgica used to carry)
$ws.Range("A30").Value = "Paises Bajos"
$ws.Range("B30").Value = 253134
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 6873

# --- Kazajistan (row 45) ---
$ws.Range("B45").Value = 109907
$ws.Range("C45").Value = 141
$ws.Range("D45").Value = 105385
$ws.Range("E45").Value = 2726

# --- Honduras (row 54) ---
$ws.Range("B54").Value = 91078
$ws.Range("C54").Value = 846
$ws.Range("D54").Value = 36341
$ws.Range("E54").Value = 52141
$ws.Range("G54").Value = 14
$ws.Range("H54").Value = 2596

# --- Venezuela (row 56) ---
$ws.Range("B56").Value = 88035
$ws.Range("D56").Value = 81626
$ws.Range("E56").Value = 5656
$ws.Range("H56").Value = 753

# --- Belice (row 153) ---
$ws.Range("B153").Value = 2937
$ws.Range("C153").Value = 51
$ws.Range("D153").Value = 1756
$ws.Range("E153").Value = 1135
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 46

# --- San Martin (Parte Holandesa) (row 172) ---
$ws.Range("B172").Value = 769
$ws.Range("C172").Value = 7
$ws.Range("D172").Value = 681
$ws.Range("E172").Value = 66
